$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("F1").Value = "height"
$ws.Range("G1").Value = "weight"

# Match the header style (bold, thin border, centered) already used by B1:E1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)   # xlPasteFormats

# Update column E values and add F/G values for rows 2-17
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667   # column E
    $ws.Cells.Item($r, 6).Value = 240                 # column F (height)
    $ws.Cells.Item($r, 7).Value = 0                   # column G (weight)
}
